$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RFR")
$ws.Activate()

# Update the target's Hour Angle (hour component) and Declination (degree component)
$ws.Range("A7").Value = 6
$ws.Range("A12").Value = 2

# Change G14's number format to an integer display ("0") instead of "0.000"
$ws.Range("G14").NumberFormat = "0"

# Move the selection/active cell to G15 (matches the saved view state)
$ws.Range("G15").Select()
